# Commit message: CO cpu fix problem about "1"
#
# The "组合逻辑译码表" (combinational-logic decode table) sheet is the
# workbook's active sheet. Several decode-condition cells only tested the
# bare "!ST" / "ST" signal, but the correct condition also needs to gate on
# "W1" ("1" in the commit message refers to this W1 signal). This script
# fixes those cells to read "!ST&W1" / "ST&W1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that read "!ST" and must become "!ST&W1"
$notStCells = @("M2", "M4", "P4", "Q4", "M5", "M11", "P11", "P14", "Q14", "M18")
foreach ($addr in $notStCells) {
    $ws.Range($addr).Value = "!ST&W1"
}

# Cells that read "ST" and must become "ST&W1"
$stCells = @("P12", "Q13", "P16", "Q16")
foreach ($addr in $stCells) {
    $ws.Range($addr).Value = "ST&W1"
}

# Keep the sheet active with the same selection as the authored edit.
$ws.Activate()
$ws.Range("L5").Select()
